$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Append a new paragraph "ura ura" after the last paragraph ("kak dela")
# ------------------------------------------------------------------
$end = $d.Content
$end.Collapse(0)          # wdCollapseEnd
$end.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(1)     # wdCollapseStart
$newRange.InsertAfter("ura ura")

# ------------------------------------------------------------------
# 2. Normal style paragraph formatting tweaks
# ------------------------------------------------------------------
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.Hyphenation = $false
$normal.ParagraphFormat.SpaceBefore = 0
$normal.ParagraphFormat.SpaceAfter = 0
$normal.ParagraphFormat.Alignment = 0
